# Data_Source_Exc.xlsx - "added some more rows"
# The sheet holds a small UPI bank-statement export. Row 4 previously held a
# stray placeholder ("my name is ali"); it is replaced with a real
# transaction row, and five more transaction rows are appended after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

# --- 1. Fix the date on row 2 (was 45566, should be 45564) -----------------
$ws.Range("A2").Value = 45564

# --- 2. Harmonise formatting for the whole data block (rows 2-9) to match
#        row 3's existing look (date cell without wrap, etc.) before filling
#        in the new cells, so every row ends up visually consistent.
$ws.Range("A3:E3").Copy()
$ws.Range("A2:E9").PasteSpecial(-4122)

# --- 3. Row 4: was the "my name is ali" placeholder, now a real row -------
$ws.Range("A4").Value = 45567
$ws.Range("B4").Value = "UPI/KAZI SANA KAZI /078513505036/Payment from Ph"
$ws.Range("C4").Value = "UPI-427637049656"
$ws.Range("D4").Value = "-" + $nbsp + "₹30.00"
$ws.Range("E4").Value = "₹2.72"

# --- 4. New row 5 -----------------------------------------------------------
$ws.Range("A5").Value = 45567
$ws.Range("B5").Value = "UPI/KAZI SANA KAZI /344059039001/Payment from Ph"
$ws.Range("C5").Value = "UPI-427620658857"
$ws.Range("D5").Value = "-" + $nbsp + "₹10.00"
$ws.Range("E5").Value = "₹32.72"

# --- 5. New row 6 -----------------------------------------------------------
$ws.Range("A6").Value = 45564
$ws.Range("B6").Value = "UPI/Kazi Shoeboddin/395535689576/Payment from Ph"
$ws.Range("C6").Value = "UPI-427376235606"
$ws.Range("D6").Value = "-" + $nbsp + "₹34.00"
$ws.Range("E6").Value = "₹117.72"

# --- 6. New row 7 -----------------------------------------------------------
$ws.Range("A7").Value = 45564
$ws.Range("B7").Value = "UPI/Shaikh Mumtaz F/038147999834/Payment from Ph"
$ws.Range("C7").Value = "UPI-427376185816"
$ws.Range("D7").Value = "-" + $nbsp + "₹166.00"
$ws.Range("E7").Value = "₹151.72"

# --- 7. New row 8 (a credit, so column D keeps the normal/black style, not
#        the red "debit" style -> copy C8's format onto D8 after filling it) -
$ws.Range("A8").Value = 45564
$ws.Range("B8").Value = "UPI/KaziShoeboddinM/463981491028/PaymentfromPhon"
$ws.Range("C8").Value = "UPI-427375877465"
$ws.Range("D8").Value = "₹200.00"
$ws.Range("E8").Value = "₹317.72"
$ws.Range("C8").Copy()
$ws.Range("D8").PasteSpecial(-4122)

# --- 8. New row 9 ------------------------------------------------------------
$ws.Range("A9").Value = 45563
$ws.Range("B9").Value = "UPI/KAZI SANA KAZI /322420393826/Payment from Ph"
$ws.Range("C9").Value = "UPI-427264687389"
$ws.Range("D9").Value = "-" + $nbsp + "₹50.00"
$ws.Range("E9").Value = "₹117.72"

# --- 9. Row heights for the new rows match the existing transaction rows ---
$ws.Rows.Item(4).RowHeight = 42.75
$ws.Rows.Item(5).RowHeight = 42.75
$ws.Rows.Item(6).RowHeight = 42.75
$ws.Rows.Item(7).RowHeight = 42.75
$ws.Rows.Item(8).RowHeight = 42.75
$ws.Rows.Item(9).RowHeight = 42.75

# --- 10. Selection ends on E9, matching where entry left off ---------------
$ws.Range("E9").Select()
